$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Title (Heading1) and the later bold "title" run reuse the same text.
# Replace=2 (wdReplaceAll) rewrites every matching occurrence in $d.Content,
# so a single call updates both the heading and the bold run below.
Replace-Text "Play Fortunes of Ali Baba Free – Slot Game Review" "Play Fortunes of Ali Baba Slot Game for Free"

# Meta description (italic run) - replace this full sentence first, since it
# contains substrings ("Den of Thieves and Free Spins feature...") that would
# otherwise collide with the shorter list-item replacements below.
Replace-Text "Read our review of Fortunes of Ali Baba slot game. Play for free and enjoy the Den of Thieves and Free Spins features, Walking Wilds, and high-quality graphics." "Read our review of Fortunes of Ali Baba and play this exciting slot game for free."

# "What we like" list
Replace-Text "Den of Thieves and Free Spins feature" "Den of Thieves Bonus feature"
Replace-Text "Walking Wilds feature increases chances of winning" "Free Spins feature"
Replace-Text "Gameplay on all devices" "Walking Wilds feature"
Replace-Text "High-quality graphics and immersive backdrop" "Gameplay on All Devices"

# "What we don't like" list
Replace-Text "Being caught by the thieves takes you back to the base game" "Potential to be caught by thieves in Den of Thieves Bonus feature"
Replace-Text "Only 8 free spins initially" "Limited number of free spins"
